$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add Thickness of Material ---
# Insert a new column before the current "Isolants" column (B), which
# shifts "Isolants" (and its data below) from column B to column C and
# frees up a blank column D.
$ws.Columns.Item(2).Insert()

# Give the two new header cells (B1, D1) the same formatting as the
# existing header cells (bold, centered, bordered) by copying the
# format from A1, then set their text.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B1").Value = "Épaisseur Matériaux (mm)"
$ws.Range("D1").Value = "Épaisseur Isolants (mm)"

# Thickness (mm) of the material "Bois" on row 2
$ws.Range("B2").Value = 0.012

# Thickness (mm) of the insulator "Laine minérale" on row 2
$ws.Range("D2").Value = 0.07000000000000001

# New row 3: a second "Laine minérale" entry with its own thickness
$ws.Range("A3").Borders.LineStyle = 0
$ws.Range("B3").Borders.LineStyle = 0
$ws.Range("C3").Value = "Laine minérale"
$ws.Range("D3").Value = 0.04
